# Update cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numbers (e.g. "0.9970", "24.667.12").
# Force the column to Text format first so Excel does not silently convert
# these values to native numbers (which would drop trailing zeros / thousand
# separators and change the cell type away from string).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '24.667.12'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.705.47'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D4').Value = '0.9970'
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').Value = '315.99'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = '0.9956'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('D7').Value = '0.3914'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = '0.4066'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').Value = '1.491'
$ws.Range('E9').Value = '  -2.14%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '53.89'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('B11').Value = 'BinanceUSD'
$ws.Range('C11').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D11').Value = '0.9946'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').Value = '0.08809'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').Value = '26.17'
$ws.Range('E13').Value = '  +10.63%  '
$ws.Range('D14').Value = '7.471'
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '8.132'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').Value = '1.700.72'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '97.88'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').Value = '0.07183'
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('E20').Value = '  +4.05%  '
$ws.Range('D21').Value = '7.313'
$ws.Range('E21').Value = '  +3.26%  '
$ws.Range('D22').Value = '0.9963'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '14.35'
$ws.Range('E23').Value = '  -2.71%  '
$ws.Range('D24').Value = '24.656.94'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '3.034'
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('D26').Value = '2.328'
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('D27').Value = '23.05'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').Value = '167.15'
$ws.Range('E28').Value = '  +2.01%  '
$ws.Range('E29').Value = '  +16.18%  '
$ws.Range('D30').Value = '144.29'
$ws.Range('E30').Value = '  +6.36%  '
$ws.Range('D31').Value = '8.422'
$ws.Range('E31').Value = '  -4.30%  '
$ws.Range('D32').Value = '1.888.81'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').Value = '0.08802'
$ws.Range('E33').Value = '  -2.16%  '
$ws.Range('E34').Value = '  +10.19%  '
$ws.Range('D35').Value = '1.060'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').Value = '7.241'
$ws.Range('E36').Value = '  -4.94%  '
$ws.Range('D37').Value = '0.03108'
$ws.Range('E37').Value = '  +5.14%  '
$ws.Range('D38').Value = '0.8744'
$ws.Range('E38').Value = '  +14.01%  '
$ws.Range('D39').Value = '0.2807'
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').Value = '10.92'
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('D41').Value = '0.09192'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').Value = '14.25'
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('D43').Value = '1.482'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').Value = '17.40'
$ws.Range('E44').Value = '  +8.72%  '
$ws.Range('D45').Value = '0.7507'
$ws.Range('E45').Value = '  +4.57%  '
$ws.Range('D46').Value = '2.683'
$ws.Range('E46').Value = '  +3.59%  '
$ws.Range('D47').Value = '4.250'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '1.397'
$ws.Range('E48').Value = '  +3.78%  '
$ws.Range('D49').Value = '0.9969'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '140.81'
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('D51').Value = '0.08259'
$ws.Range('E51').Value = '  +3.53%  '
